# Development status.xlsx - add a new data row (week of 03/13) and extend
# the two line charts so they plot through the new row; also reposition
# the charts (they were dragged to new anchor cells in the source edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Append the new data row (row 6) under the existing table (A1:E5).
# ---------------------------------------------------------------------
$ws.Range("A6").Value2 = 42076
$ws.Range("B6").Value2 = 21
$ws.Range("C6").Value2 = 89
$ws.Range("D6").Value2 = 7188
$ws.Range("E6").Value2 = 10377

# Carry the date-number formatting from A5 down onto the new A6 cell
# (copy formats only, so the values just written above are untouched).
$null = $ws.Range("A5").Copy()
$null = $ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Move the selection the same way the original edit left it.
$null = $ws.Range("E19").Select()

# ---------------------------------------------------------------------
# 2) Extend both line charts' series so they include the new row.
# ---------------------------------------------------------------------
$cos = $ws.ChartObjects()

# Chart 1: "Packages & classes" -> series B (Packages) and C (Classes)
$chart1 = $cos.Item(1).Chart
$series1 = $chart1.SeriesCollection()
$series1.Item(1).Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$6,Sheet1!`$B`$2:`$B`$6,1)"
$series1.Item(2).Formula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$A`$6,Sheet1!`$C`$2:`$C`$6,2)"

# Chart 2: "Source code lines" -> series D (Method lines) and E (Total lines)
$chart2 = $cos.Item(2).Chart
$series2 = $chart2.SeriesCollection()
$series2.Item(1).Formula = "=SERIES(Sheet1!`$D`$1,Sheet1!`$A`$2:`$A`$6,Sheet1!`$D`$2:`$D`$6,1)"
$series2.Item(2).Formula = "=SERIES(Sheet1!`$E`$1,Sheet1!`$A`$2:`$A`$6,Sheet1!`$E`$2:`$E`$6,2)"

# ---------------------------------------------------------------------
# 3) Reposition the two chart objects on the sheet (they were dragged to
#    new anchor cells to make room for the taller data table).
# ---------------------------------------------------------------------
$co1 = $cos.Item(1)
$co1.Left = 333.1875
$co1.Top = 17.25
$co1.Width = 391.0625
$co1.Height = 216.0

$co2 = $cos.Item(2)
$co2.Left = 335.06246062992125
$co2.Top = 249.0
$co2.Width = 391.0625
$co2.Height = 216.0

Write-Host "Row 6 appended; charts rescoped to A2:A6/B2:B6/.../E2:E6 and repositioned."
